function Get-ParaIndexByText($doc, $pattern) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text -like $pattern) {
            return $i
        }
    }
    return -1
}

$d = $word.ActiveDocument

# --- Change 1: after the "maven 3.9.9" paragraph, add a new paragraph
#     "Docker version 28.0.1," that carries the _GoBack bookmark.
#     (Word keeps only one "_GoBack" bookmark at a time; adding a new one
#     here automatically removes it from its old location further down.)
$mavenIdx = Get-ParaIndexByText $d "maven 3.9.9*"
$mavenPara = $d.Paragraphs.Item($mavenIdx)
$mavenPara.Range.InsertParagraphAfter()

$dockerPara = $d.Paragraphs.Item($mavenIdx + 1)
$dockerPara.Range.Text = "Docker version 28.0.1,"
$dockerPara = $d.Paragraphs.Item($mavenIdx + 1)
$dockerTextRange = $d.Range($dockerPara.Range.Start, $dockerPara.Range.End - 1)
$d.Bookmarks.Add("_GoBack", $dockerTextRange)

# --- Change 2: move <w:lastRenderedPageBreak/> from the run containing
#     "cd history" to the run containing "Frontend" (the bold heading just
#     before it), keeping all other paragraph/run formatting intact.
$frontendIdx = Get-ParaIndexByText $d "Frontend`r"
$frontendPara = $d.Paragraphs.Item($frontendIdx)
$cdHistoryPara = $d.Paragraphs.Item($frontendIdx + 1)

$frontendRunRange = $d.Range($frontendPara.Range.Start, $frontendPara.Range.End - 1)
$frontendRunRange.InsertXML("<pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`"><pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`"><pkg:xmlData><w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:body><w:p><w:r w:rsidRPr=`"00213D94`"><w:rPr><w:b/></w:rPr><w:lastRenderedPageBreak/><w:t>Frontend</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>")

$cdHistoryRunRange = $d.Range($cdHistoryPara.Range.Start, $cdHistoryPara.Range.End - 1)
$cdHistoryRunRange.InsertXML("<pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`"><pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`"><pkg:xmlData><w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:body><w:p><w:r><w:t>cd history</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>")
